$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = "Class A"
$ws.Range("D3").Value = "Class A*"

$ws.Range("D3").Select()
